$d = $word.ActiveDocument

# Locate the paragraph that starts the "Em maio de 2011..." text.
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Em maio de 2011 iniciei")) {
        $idx = $i
        break
    }
}

# Merge the following empty paragraph into this one (removes the blank
# paragraph that used to separate the two sentences) by deleting this
# paragraph's end-of-paragraph mark.
$p = $d.Paragraphs.Item($idx)
$rngEnd = $d.Range($p.Range.End - 1, $p.Range.End)
$rngEnd.Delete()

# Merge the next paragraph ("Neste período de 2011 iniciei...") into it
# as well, joining the two sentences into a single paragraph.
$p = $d.Paragraphs.Item($idx)
$rngEnd2 = $d.Range($p.Range.End - 1, $p.Range.End)
$rngEnd2.Delete()

# Insert the separating space between the two former sentences, and
# change "iniciei" to "havia iniciado" right before "a faculdade".
$d.Content.Find.Execute("dispensado.Neste", $true, $false, $false, $false, $false, $true, 1, $false, "dispensado. Neste", 2) | Out-Null
$d.Content.Find.Execute("2011 iniciei a faculdade", $true, $false, $false, $false, $false, $true, 1, $false, "2011 havia iniciado a faculdade", 2) | Out-Null
